$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $savedStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $savedStyle
}

Set-TextValue "D2" "27.644.27"
Set-TextValue "E2" "  -0.50%  "
Set-TextValue "D3" "1.617.69"
Set-TextValue "E3" "  -0.62%  "
Set-TextValue "D4" "0.990"
Set-TextValue "E4" "  -0.51%  "
Set-TextValue "D5" "209.34"
Set-TextValue "E5" "  -0.83%  "
Set-TextValue "D6" "0.516"
Set-TextValue "E6" "  -1.02%  "
Set-TextValue "E7" "  -0.48%  "
Set-TextValue "D8" "23.11"
Set-TextValue "E8" "  -0.52%  "
Set-TextValue "E9" "  -0.95%  "
Set-TextValue "E10" "  -1.21%  "
Set-TextValue "D11" "0.0875"
Set-TextValue "E11" "  -0.53%  "
Set-TextValue "D12" "1.846.47"
Set-TextValue "E12" "  -0.78%  "
Set-TextValue "D13" "1.627.81"
Set-TextValue "E13" "  -0.14%  "
Set-TextValue "E14" "  -1.42%  "
Set-TextValue "D15" "0.556"
Set-TextValue "E15" "  -1.17%  "
Set-TextValue "D16" "64.57"
Set-TextValue "E16" "  -0.89%  "
Set-TextValue "D17" "27.654.41"
Set-TextValue "E17" "  -0.59%  "
Set-TextValue "D18" "227.60"
Set-TextValue "E18" "  -0.98%  "
Set-TextValue "D19" "7.65"
Set-TextValue "E19" "  +1.85%  "
Set-TextValue "E20" "  -0.98%  "
Set-TextValue "D21" "0.990"
Set-TextValue "E21" "  -0.50%  "
Set-TextValue "D22" "4.31"
Set-TextValue "E22" "  -1.06%  "
Set-TextValue "D23" "10.07"
Set-TextValue "E23" "  -2.17%  "
Set-TextValue "D24" "2.03"
Set-TextValue "E24" "  -0.90%  "
Set-TextValue "D25" "154.26"
Set-TextValue "E25" "  +0.17%  "
Set-TextValue "D26" "6.88"
Set-TextValue "E26" "  -0.82%  "
Set-TextValue "E27" "  -0.52%  "
Set-TextValue "D28" "15.41"
Set-TextValue "E28" "  -1.21%  "
Set-TextValue "D29" "0.990"
Set-TextValue "E29" "  -0.54%  "
Set-TextValue "E30" "  -0.63%  "
Set-TextValue "E31" "  -0.57%  "
Set-TextValue "E32" "  -1.05%  "
Set-TextValue "E33" "  -0.04%  "
Set-TextValue "D34" "1.391.72"
Set-TextValue "E34" "  -0.74%  "
Set-TextValue "D35" "1.59"
Set-TextValue "E35" "  +1.60%  "
Set-TextValue "D36" "0.995"
Set-TextValue "E36" "  -1.66%  "
Set-TextValue "E37" "  -1.40%  "
Set-TextValue "E38" "  +0.39%  "
Set-TextValue "D39" "0.555"
Set-TextValue "E39" "  -0.87%  "
Set-TextValue "D40" "0.842"
Set-TextValue "E40" "  -2.85%  "
Set-TextValue "E41" "  -1.16%  "
Set-TextValue "E42" "  -0.55%  "
Set-TextValue "E43" "  -0.21%  "
Set-TextValue "D44" "65.51"
Set-TextValue "E44" "  -1.58%  "
Set-TextValue "D45" "5.36"
Set-TextValue "E45" "  -2.59%  "
Set-TextValue "D46" "1.755.78"
Set-TextValue "E46" "  -1.13%  "
Set-TextValue "E47" "  -3.35%  "
Set-TextValue "D48" "87.59"
Set-TextValue "E48" "  -0.06%  "
Set-TextValue "E49" "  +1.41%  "
Set-TextValue "E50" "  -0.56%  "
Set-TextValue "D51" "7.55"
Set-TextValue "E51" "  +1.34%  "
